$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") was bumped by one day (2023-09-10 -> 2023-09-11,
# serial 45179 -> 45180) for every data row (rows 2-221).
$ws.Range("C2:C221").Value = 45180
